$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update iteration data rows (A2:D13 existing, A14:D16 new) with text-valued numbers
# matching the inlineStr encoding used by the source table.
$data = @{
    2 = @("0", "2.0", "-9.85203026391962", "1.0000005")
    3 = @("1", "11.8520302639196", "-12.4554640184409", "9.85203026391962")
    4 = @("2", "24.3074942823605", "-5.02799868400258", "12.4554640184409")
    5 = @("3", "29.3354929663631", "-1.31609406230989", "5.02799868400258")
    6 = @("4", "30.651587028673", "-0.307204052149675", "1.31609406230989")
    7 = @("5", "30.9587910808226", "-0.0698079204520248", "0.307204052149675")
    8 = @("6", "31.0285990012747", "-0.0157662921150674", "0.0698079204520248")
    9 = @("7", "31.0443652933897", "-0.0035559457626313", "0.0157662921150674")
    10 = @("8", "31.0479212391524", "-0.0008017620788542", "0.0035559457626313")
    11 = @("9", "31.0487230012312", "-0.0001807612833495", "0.0008017620788542")
    12 = @("10", "31.0489037625146", "-4.07528934509571e-05", "0.0001807612833495")
    13 = @("11", "31.048944515408", "-9.18776614966532e-06", "4.07528934509571e-05")
    14 = @("12", "31.0489537031742", "-2.0713861452748e-06", "9.18776614966532e-06")
    15 = @("13", "31.0489557745603", "-4.66994883652205e-07", "2.0713861452748e-06")
    16 = @("14", "31.0489562415552", "-1.05284190254906e-07", "4.66994883652205e-07")
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($col = 1; $col -le 4; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $vals[$col - 1]
        $cell.Style = "Normal"
    }
}

Write-Host "Used range:" $ws.UsedRange.Address()
